# modulate factor of hough circle

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the C column ("modulated" hough-circle factor) values for rows 3-32
$values = [ordered]@{
    3  = 43
    4  = 59
    5  = 41
    6  = 106
    7  = 39
    8  = 23
    9  = 198
    10 = 176
    11 = 187
    12 = 136
    13 = 92
    14 = 206
    15 = 152
    16 = 3
    17 = 32
    18 = 1
    19 = 170
    20 = 88
    21 = 91
    22 = 101
    23 = 72
    24 = 78
    25 = 11
    26 = 24
    27 = 195
    28 = 8
    29 = 213
    30 = 186
    31 = 181
    32 = 186
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

# Move the active selection on the sheet to G12
$ws.Activate()
$ws.Range("G12").Select()

# Reposition / resize the workbook window to match the saved view state
$win = $excel.ActiveWindow
$win.Left = 3060
$win.Top = 2625
$win.Width = 12345
$win.Height = 11385
